{"js": "// Word JS API (Office.js) edit script.\n//\n// Change being applied (see commit message \"notes of Listening of IELTS 7\n// Test 1\"):\n//   1. The \"on a roll\" idiom paragraph loses the stray paragraph-mark\n//      formatting (<w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr> that used\n//      to live in its <w:pPr>), and it no longer carries the `_GoBack`\n//      bookmark at its end.\n//   2. A brand-new idiom note - \"come to think of it.\" - is appended as a\n//      new last paragraph (it reuses what used to be the trailing empty\n//      paragraph), and the `_GoBack` bookmark now sits at the end of that\n//      paragraph instead.\n//\n// The \"strike/touch a chord with somebody\" paragraph in between is\n// untouched content-wise.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two paragraphs we need by their (stable, unique) text rather\n// than a hard-coded index.\nlet onARollIndex = -1;\nlet lastIndex = paragraphs.items.length - 1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"on a roll\") === 0) {\n    onARollIndex = i;\n    break;\n  }\n}\nif (onARollIndex === -1) {\n  throw new Error('Could not find the \"on a roll\" paragraph.');\n}\n\nconst onARollParagraph = paragraphs.items[onARollIndex];\nconst lastParagraph = paragraphs.items[lastIndex];\n\n// Helper: wrap a <w:body> fragment in the minimal package envelope that\n// Word's insertOoxml expects.\nfunction wrapOoxml(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    bodyXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// 1) \"on a roll\" paragraph: drop the pPr/rPr eastAsia hint and the\n//    trailing _GoBack bookmark; keep every run exactly as it was.\nconst onARollXml =\n  \"<w:p>\" +\n  '<w:pPr><w:jc w:val=\"left\"/></w:pPr>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:b/></w:rPr><w:t>on a roll</w:t></w:r>' +\n  \"<w:r><w:br/></w:r>\" +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>to be having a successful or lucky period.</w:t></w:r>' +\n  \"<w:r><w:br/></w:r>\" +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\">Pippa won five games in a row and it was </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>evident</w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> that she was on a roll.</w:t></w:r>' +\n  \"</w:p>\";\n\nonARollParagraph.insertOoxml(wrapOoxml(onARollXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Trailing (previously empty) paragraph: add the new \"come to think of\n//    it.\" note and move the _GoBack bookmark here.\nconst lastParagraphXml =\n  \"<w:p>\" +\n  '<w:pPr><w:jc w:val=\"left\"/></w:pPr>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>come to think of it.</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  \"</w:p>\";\n\nlastParagraph.insertOoxml(wrapOoxml(lastParagraphXml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) edit script.\n#\n# Change being applied (see commit message \"notes of Listening of IELTS 7\n# Test 1\"):\n#   1. The \"on a roll\" idiom paragraph loses the stray paragraph-mark\n#      formatting (<w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr> that used\n#      to live in its <w:pPr>), and it no longer carries the `_GoBack`\n#      bookmark at its end.\n#   2. A brand-new idiom note - \"come to think of it.\" - is appended as a\n#      new last paragraph (it reuses what used to be the trailing empty\n#      paragraph), and the `_GoBack` bookmark now sits at the end of that\n#      paragraph instead.\n#\n# The \"strike/touch a chord with somebody\" paragraph in between is\n# untouched content-wise.\n\n$d = $word.ActiveDocument\n\n$wordNs = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"\n\n# Locate the \"on a roll\" paragraph and the trailing (last, empty) paragraph\n# by content rather than a hard-coded index.\n$onARoll = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"on a roll\")) {\n        $onARoll = $p\n        break\n    }\n}\nif ($onARoll -eq $null) {\n    throw \"Could not find the 'on a roll' paragraph.\"\n}\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n\n# 1) \"on a roll\" paragraph: drop the pPr/rPr eastAsia hint and the\n#    trailing _GoBack bookmark; keep every run exactly as it was.\n$onARollXml = (\n    '<w:p xmlns:w=\"' + $wordNs + '\">' +\n    '<w:pPr><w:jc w:val=\"left\"/></w:pPr>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:b/></w:rPr><w:t>on a roll</w:t></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>to be having a successful or lucky period.</w:t></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\">Pippa won five games in a row and it was </w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>evident</w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> that she was on a roll.</w:t></w:r>' +\n    '</w:p>'\n)\n[void]$onARoll.Range.InsertXML($onARollXml)\n\n# 2) Trailing (previously empty) paragraph: add the new \"come to think of\n#    it.\" note and move the _GoBack bookmark here.\n$lastParagraphXml = (\n    '<w:p xmlns:w=\"' + $wordNs + '\">' +\n    '<w:pPr><w:jc w:val=\"left\"/></w:pPr>' +\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>come to think of it.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\n)\n[void]$lastParagraph.Range.InsertXML($lastParagraphXml)\n"}
